$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''93.618.54'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.14%  '

$ws.Range("D3").Value = '''3.424.60'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.98%  '

$ws.Range("E4").Value = '  -0.15%  '

$ws.Range("D5").Value = '''234.29'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.53%  '

$ws.Range("D6").Value = '''621.88'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.19%  '

$ws.Range("E7").Value = '  -1.72%  '

$ws.Range("E8").Value = '  -0.09%  '

$ws.Range("E9").Value = '  -0.02%  '

$ws.Range("D10").Value = '''0.977'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.31%  '

$ws.Range("D11").Value = '''3.421.00'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.90%  '

$ws.Range("D12").Value = '''43.05'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.98%  '

$ws.Range("D13").Value = '''0.200'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.92%  '

$ws.Range("D14").Value = '''6.30'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.90%  '

$ws.Range("D15").Value = '''93.420.02'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.20%  '

$ws.Range("D16").Value = '''4.069.13'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.95%  '

$ws.Range("E17").Value = '  +0.24%  '

$ws.Range("E18").Value = '  +0.44%  '

$ws.Range("D19").Value = '''3.420.52'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.81%  '

$ws.Range("D20").Value = '''18.09'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.27%  '

$ws.Range("D21").Value = '''11.67'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.65%  '

$ws.Range("E22").Value = '  +5.64%  '

$ws.Range("D23").Value = '''503.32'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.81%  '

$ws.Range("E24").Value = '  -2.68%  '

$ws.Range("D25").Value = '''6.66'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.28%  '

$ws.Range("E26").Value = '  -2.05%  '

$ws.Range("D27").Value = '''95.16'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.00%  '

$ws.Range("D28").Value = '''12.01'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.32%  '

$ws.Range("D29").Value = '''3.606.16'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.76%  '

$ws.Range("D30").Value = '''11.43'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.65%  '

$ws.Range("E31").Value = '  +0.00%  '

$ws.Range("D32").Value = '''0.139'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.07%  '

$ws.Range("E33").Value = '  +2.37%  '

$ws.Range("D34").Value = '''0.996'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.62%  '

$ws.Range("E35").Value = '  -0.27%  '

$ws.Range("D36").Value = '''29.96'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.02%  '

$ws.Range("D37").Value = '''0.552'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.53%  '

$ws.Range("D38").Value = '''559.18'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.35%  '

$ws.Range("E39").Value = '  -0.62%  '

$ws.Range("E40").Value = '  -0.54%  '

$ws.Range("E41").Value = '  -0.01%  '

$ws.Range("D42").Value = '''0.150'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.66%  '

$ws.Range("D43").Value = '''0.915'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.03%  '

$ws.Range("E44").Value = '  +2.66%  '

$ws.Range("D45").Value = '''23.68'
$ws.Range("D45").Style = "Normal"

$ws.Range("D46").Value = '''3.68'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.19%  '

$ws.Range("D47").Value = '''0.0412'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.29%  '

$ws.Range("D48").Value = '''5.50'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.29%  '

$ws.Range("D49").Value = '''53.64'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.58%  '

$ws.Range("E50").Value = '  -1.03%  '

$ws.Range("E51").Value = '  +2.34%  '
